$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The viability factor is introduced right after the existing "mX0" column
# (H=cX0, I=mX0). Two new columns are inserted there: the old "mX0" formula
# column is relabelled "mX0_ohne_Viab_f" (stays in place, column I), a new
# "viab_f" column is added (J), and a new "mX0" column recomputes the
# viability-corrected biomass as I*J (K). Everything that used to live at
# J onward shifts two columns to the right.
$ws.Range("J1:K1").EntireColumn.Insert() | Out-Null

# Give the two freshly inserted columns the same width as the neighbouring
# H:I columns (raw column width ~11 units <-> COM ColumnWidth 10.17).
$ws.Range("J1:K1").ColumnWidth = 10.17

# --- Header row (row 1) ---
$ws.Range("I1").Value = "mX0_ohne_Viab_f"
$ws.Range("J1").Value = "viab_f"
$ws.Range("K1").Value = "mX0"

# --- Unit row (row 2) ---
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "g"

# --- Data rows (3 and 4): new viability factor values + recomputed mX0 ---
$ws.Range("J3").NumberFormat = "0.000"
$ws.Range("J3").Value = 0.0459
$ws.Range("K3").NumberFormat = "0.000"
$ws.Range("K3").Formula = "=I3*J3"

$ws.Range("J4").NumberFormat = "0.000"
$ws.Range("J4").Value = 0.3244
$ws.Range("K4").NumberFormat = "0.000"
$ws.Range("K4").Formula = "=I4*J4"

# Restore the selection to where the author last left the cursor.
$ws.Range("L9").Select() | Out-Null
